$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the MKO10 sample row (originally row 33). Select it first so the
# resulting selection matches where the user had been working.
$ws.Rows.Item(33).Select() | Out-Null
$ws.Rows.Item(33).Delete() | Out-Null

# Delete the MYWT-6 sample row (originally the very last row, 54;
# after the previous deletion it has become row 53).
$ws.Rows.Item(53).Delete() | Out-Null
